# Record a new day of progress (04-Nov-2025): every training row's
# "PERIOD TO EXPIRE" (column H) ticks down by one day, and the
# "LAST UPDATE" (column I) date moves from 03-Nov-2025 to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Scratch cell used to inject the new date as plain text (not an Excel
# date serial) while keeping the destination cells' existing style.
# Formatting the scratch cell as Text and then pasting only the value
# into the target cell avoids Excel's automatic "looks like a date"
# conversion while leaving the target cell's style/borders untouched.
$scratch = $ws.Cells.Item(100, 1)
$scratch.NumberFormat = "@"

for ($row = 3; $row -le 25; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H: PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # column I: LAST UPDATE

    if ($iCell.Value2 -eq "03-Nov-2025") {
        $hCell.Value = $hCell.Value2 - 1

        $scratch.Value = "04-Nov-2025"
        $scratch.Copy()
        $iCell.PasteSpecial(-4163)   # xlPasteValues
    }
}

$scratch.Clear()
$excel.CutCopyMode = $false
